$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.386.93'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '1.565.87'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.33%  '
$ws.Range('E6').Value = '  -0.93%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('E9').Value = '  -1.83%  '
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0867'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').Value = '1.788.57'
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('D13').Value = '1.572.43'
$ws.Range('E13').Value = '  -1.51%  '
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('E15').Value = '  -2.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.69%  '
$ws.Range('D17').Value = '27.390.47'
$ws.Range('E17').Value = '  -0.97%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.0₃0688'
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '211.97'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.67%  '
$ws.Range('E20').Value = '  -1.12%  '
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.11'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.51'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.35%  '
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.54'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.72'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('E29').Value = '  -1.94%  '
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('E31').Value = '  +1.11%  '
$ws.Range('E32').Value = '  -0.92%  '
$ws.Range('D33').Value = '1.364.46'
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('E34').Value = '  +0.19%  '
$ws.Range('E35').Value = '  +1.41%  '
$ws.Range('E36').Value = '  +0.39%  '
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('E38').Value = '  +0.69%  '
$ws.Range('E39').Value = '  -0.97%  '
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.973'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.42%  '
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '64.03'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.29'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.55%  '
$ws.Range('E46').Value = '  -1.36%  '
$ws.Range('D47').Value = '1.701.16'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.60'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Value = '0.0₇0993'
$ws.Range('E49').Value = '  +3.18%  '
$ws.Range('E51').Value = '  -0.64%  '
